$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.86
$ws.Range("H2").Value = 2.48
$ws.Range("J2").Value = 3.3
$ws.Range("P2").Value = 1.87
$ws.Range("AA2").Value = 980
$ws.Range("AE2").Value = 980
$ws.Range("AI2").Value = 980
$ws.Range("AJ2").Value = 980
$ws.Range("AK2").Value = 980
$ws.Range("AL2").Value = 980
$ws.Range("AN2").Value = 980
$ws.Range("N3").Value = 5.5
$ws.Range("O3").Value = 1.2
$ws.Range("P3").Value = 2.48
$ws.Range("Q3").Value = 1.64
$ws.Range("R3").Value = 1.6
$ws.Range("S3").Value = 2.56
$ws.Range("T3").Value = 1.62
$ws.Range("U3").Value = 2.54
$ws.Range("AA3").Value = 90
$ws.Range("AB3").Value = 12.5
$ws.Range("AC3").Value = 9.4
$ws.Range("AG3").Value = 10
$ws.Range("AH3").Value = 16
$ws.Range("AK3").Value = 17
$ws.Range("AO3").Value = 36
$ws.Range("H4").Value = 1.32
$ws.Range("I4").Value = 1.35
$ws.Range("J4").Value = 5.6
$ws.Range("P4").Value = 2.34
$ws.Range("V4").Value = 3.85
$ws.Range("F5").Value = 1.83
$ws.Range("K5").Value = 4.5
$ws.Range("Y5").Value = 980
$ws.Range("AD5").Value = 980
$ws.Range("AH5").Value = 980
$ws.Range("AJ5").Value = 980
$ws.Range("AK5").Value = 980
$ws.Range("AL5").Value = 980
$ws.Range("F6").Value = 1.38
$ws.Range("G6").Value = 1.39
$ws.Range("K6").Value = 6
$ws.Range("S6").Value = 2.22
$ws.Range("V6").Value = 1.11
$ws.Range("W6").Value = 3.55
$ws.Range("AG6").Value = 10
$ws.Range("AK6").Value = 13
$ws.Range("L7").Value = 1.3
$ws.Range("N7").Value = 5.1
$ws.Range("R7").Value = 1.55
$ws.Range("X7").Value = 19
$ws.Range("AF7").Value = 28
$ws.Range("F8").Value = 1.7
$ws.Range("G8").Value = 1.71
$ws.Range("K8").Value = 4
$ws.Range("V8").Value = 1.18
$ws.Range("AL8").Value = 38
$ws.Range("F12").Value = 2.02
$ws.Range("H12").Value = 3.9
$ws.Range("I12").Value = 4.4
$ws.Range("J12").Value = 3.55
$ws.Range("Q12").Value = 1.9
